# Change sample process error rates
# - Reorder/rewrite the "Settings" sheet rows so ProcessRetries now sits
#   right after AppCredential (with TransactionQueue following it), and
#   ErrorsFolder/TempFolder/InputFile are pushed further down.
# - Reword the ProcessRetries and ErrorsFolder descriptions.
# - Update the active sheet/selection bookmarks to match the new state
#   (Xtras becomes the active tab; Settings/Errors selections move too).

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")

# Wipe the block that holds ErrorsFolder / TempFolder / InputFile /
# TransactionQueue / ProcessRetries plus the two stray styled-but-empty
# rows right below them (rows 17-18) so nothing is left behind once the
# content is rewritten in its new order.
$wsSettings.Range("A9:C18").Clear()

$wsSettings.Range("A10").Value = "ProcessRetries"
$wsSettings.Range("B10").Value = 2
$wsSettings.Range("C10").Value = "How many times to retry transactions in case of issues"

$wsSettings.Range("A11").Value = "TransactionQueue"
$wsSettings.Range("B11").Value = "RFW-ChorePile"
$wsSettings.Range("C11").Value = "Transactions queue in Orchestrator"

$wsSettings.Range("A13").Value = "ErrorsFolder"
$wsSettings.Range("B13").Value = "Errors"
$wsSettings.Range("C13").Value = "For system failures, logging screenshots and stack traces"

$wsSettings.Range("A14").Value = "TempFolder"
$wsSettings.Range("B14").Value = "Temp"
$wsSettings.Range("C14").Value = "Where files are kept locally while processing"

$wsSettings.Range("A16").Value = "InputFile"
$wsSettings.Range("B16").Value = "TestData\Sample1.xlsx"
$wsSettings.Range("C16").Value = "Sample data for running the template"

# Move the Settings selection, then walk through Errors and finally Xtras
# so Xtras ends up as the active tab (matches activeTab="3" / tabSelected).
$null = $wsSettings.Range("A14").Select()

$wsErrors = $wb.Worksheets.Item("Errors")
$null = $wsErrors.Range("B2").Select()

$wsXtras = $wb.Worksheets.Item("Xtras")
$null = $wsXtras.Range("B9").Select()
